$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows continue the report from 2021-09-21 (serial 44460) through
# 2021-12-08 (serial 44538), appended below the existing data (rows 386-464).
# Column A keeps the same date style (format/border/font) as the preceding
# data rows, so replicate it down before writing values.
$ws.Range("A385").Copy($ws.Range("A386:A464"))

$data = @(
  @(44460,0,0,0),
  @(44461,0,0,0),
  @(44462,0,0,0),
  @(44463,3,3,71.47962830593281),
  @(44464,0,3,71.47962830593281),
  @(44465,0,3,71.47962830593281),
  @(44466,2,5,119.1327138432213),
  @(44467,0,5,119.1327138432213),
  @(44468,0,5,119.1327138432213),
  @(44469,0,5,119.1327138432213),
  @(44470,0,2,47.65308553728854),
  @(44471,0,2,47.65308553728854),
  @(44472,0,2,47.65308553728854),
  @(44473,0,0,0),
  @(44474,0,0,0),
  @(44475,1,1,23.82654276864427),
  @(44476,1,2,47.65308553728854),
  @(44477,0,2,47.65308553728854),
  @(44478,0,2,47.65308553728854),
  @(44479,1,3,71.47962830593281),
  @(44480,1,4,95.30617107457708),
  @(44481,0,4,95.30617107457708),
  @(44482,0,3,71.47962830593281),
  @(44483,0,2,47.65308553728854),
  @(44484,0,2,47.65308553728854),
  @(44485,0,2,47.65308553728854),
  @(44486,0,1,23.82654276864427),
  @(44487,0,0,0),
  @(44488,0,0,0),
  @(44489,0,0,0),
  @(44490,0,0,0),
  @(44491,0,0,0),
  @(44492,0,0,0),
  @(44493,0,0,0),
  @(44494,0,0,0),
  @(44495,0,0,0),
  @(44496,0,0,0),
  @(44497,0,0,0),
  @(44498,0,0,0),
  @(44499,0,0,0),
  @(44500,0,0,0),
  @(44501,0,0,0),
  @(44502,0,0,0),
  @(44503,0,0,0),
  @(44504,0,0,0),
  @(44505,1,1,23.82654276864427),
  @(44506,0,1,23.82654276864427),
  @(44507,0,1,23.82654276864427),
  @(44508,0,1,23.82654276864427),
  @(44509,0,1,23.82654276864427),
  @(44510,0,1,23.82654276864427),
  @(44511,0,1,23.82654276864427),
  @(44512,0,0,0),
  @(44513,0,0,0),
  @(44514,0,0,0),
  @(44515,0,0,0),
  @(44516,1,1,23.82654276864427),
  @(44517,0,1,23.82654276864427),
  @(44518,0,1,23.82654276864427),
  @(44519,0,1,23.82654276864427),
  @(44520,0,1,23.82654276864427),
  @(44521,0,1,23.82654276864427),
  @(44522,0,1,23.82654276864427),
  @(44523,0,0,0),
  @(44524,4,4,95.30617107457708),
  @(44525,0,4,95.30617107457708),
  @(44526,0,4,95.30617107457708),
  @(44527,0,4,95.30617107457708),
  @(44528,0,4,95.30617107457708),
  @(44529,0,4,95.30617107457708),
  @(44530,0,4,95.30617107457708),
  @(44531,0,0,0),
  @(44532,2,2,47.65308553728854),
  @(44533,0,2,47.65308553728854),
  @(44534,0,2,47.65308553728854),
  @(44535,0,2,47.65308553728854),
  @(44536,1,3,71.47962830593281),
  @(44537,2,5,119.1327138432213),
  @(44538,0,5,119.1327138432213)
)

$startRow = 386
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
